$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1): P1 = 14, Q1 = 15, matching the existing bold/bordered style
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$src = $ws.Range("O1")
$dst = $ws.Range("P1:Q1")
$dst.Font.Bold = $src.Font.Bold
$dst.HorizontalAlignment = $src.HorizontalAlignment
$dst.VerticalAlignment = $src.VerticalAlignment
$dst.Borders.LineStyle = $src.Borders.LineStyle

# For data rows 2-25:
#  - flip values in columns I, K, M, O (I:1->2, K:2->1, M:1->2, O:2->1)
#  - add new columns P = 2, Q = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 15).Value = 1  # O
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
